$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Professionalism")

# --- Clear rows that fully disappear in the new layout ---
$ws.Range("A11").Clear()
$ws.Range("A15").Clear()

# --- Row 7: label text unchanged, but its cell style moves from the old
#     "22" xf onto the new variant; reuse A7's own current formatting as the
#     seed style (visually identical: Arial 12, left, wrap) for all the
#     list-item cells below. ---
$ws.Range("A7").Copy()
foreach ($addr in @("A8","A9","A10","A13","A14","A16","A17")) {
    $ws.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# --- Row 8 previously had an explicit row height (30); the new layout uses
#     the default row height, so clear that override. ---
$ws.Rows.Item(8).AutoFit()

# --- Fill in the new text content (values) ---
$ws.Range("A7").Value = "Work Process Enumeration"
$ws.Range("E7").Value = "End of Day"
$ws.Range("A8").Value = "Release Check List"
$ws.Range("A9").Value = "Version Control Check List"
$ws.Range("A10").Value = "Error Mitigation Stratagy"
$ws.Range("A12").Value = "Work Day Enumeration"
$ws.Range("A13").Value = "Task List Template"
$ws.Range("A14").Value = "Responsibility Guidelines"
$ws.Range("A16").Value = "Review Check List"
$ws.Range("A17").Value = "Conflict Guidelines"
$ws.Range("A19").Value = "Project Scoping"

# --- Update the view's active cell/selection to match the new extent ---
$ws.Activate()
$ws.Range("A7:A19").Select()
